$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row order / swapped fixtures for rows 12/13 and 15/16 ---
# (source data had these two match-pairs swapped; restore correct pairing)
$ws.Range("F12").Value = "Moghreb Tetouan"
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = "Olympique de Safi"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 3.02
$ws.Range("K12").Value = "01/09/2023 16:42"
$ws.Range("L12").Value = 3.59
$ws.Range("M12").Value = "02/09/2023 21:20"
$ws.Range("N12").Value = 2.8
$ws.Range("O12").Value = "01/09/2023 16:42"
$ws.Range("P12").Value = 2.75
$ws.Range("Q12").Value = "02/09/2023 21:20"
$ws.Range("R12").Value = 2.4
$ws.Range("S12").Value = "01/09/2023 16:42"
$ws.Range("T12").Value = 2.33
$ws.Range("U12").Value = "02/09/2023 21:20"
$ws.Range("V12").Value = "https://www.betexplorer.com/football/morocco/botola-pro/moghreb-tetouan-olympique-de-safi/Q5Wd8U5n/"

$ws.Range("F13").Value = "Mouloudia Oujda"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = "Union Touarga"
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 2.25
$ws.Range("K13").Value = "01/09/2023 16:42"
$ws.Range("L13").Value = 2.83
$ws.Range("M13").Value = "02/09/2023 21:14"
$ws.Range("N13").Value = 3.04
$ws.Range("O13").Value = "01/09/2023 16:42"
$ws.Range("P13").Value = 2.52
$ws.Range("Q13").Value = "02/09/2023 21:14"
$ws.Range("R13").Value = 3.02
$ws.Range("S13").Value = "01/09/2023 16:42"
$ws.Range("T13").Value = 3.11
$ws.Range("U13").Value = "02/09/2023 21:14"
$ws.Range("V13").Value = "https://www.betexplorer.com/football/morocco/botola-pro/mouloudia-oujda-union-touarga/4dV07lLh/"

$ws.Range("F15").Value = "FUS Rabat"
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = "Maghreb Fez"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1.73
$ws.Range("K15").Value = "02/09/2023 21:12"
$ws.Range("L15").Value = 1.85
$ws.Range("M15").Value = "03/09/2023 20:10"
$ws.Range("N15").Value = 3.23
$ws.Range("O15").Value = "02/09/2023 21:12"
$ws.Range("P15").Value = 2.83
$ws.Range("Q15").Value = "03/09/2023 20:10"
$ws.Range("R15").Value = 5.25
$ws.Range("S15").Value = "02/09/2023 21:12"
$ws.Range("T15").Value = 5.19
$ws.Range("U15").Value = "03/09/2023 20:10"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/morocco/botola-pro/fus-rabat-maghreb-fez/Obo25Sk5/"

$ws.Range("F16").Value = "Berkane"
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = "IR Tanger"
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1.98
$ws.Range("K16").Value = "02/09/2023 21:12"
$ws.Range("L16").Value = 1.56
$ws.Range("M16").Value = "03/09/2023 20:08"
$ws.Range("N16").Value = 3.11
$ws.Range("O16").Value = "02/09/2023 21:12"
$ws.Range("P16").Value = 3.59
$ws.Range("Q16").Value = "03/09/2023 20:08"
$ws.Range("R16").Value = 3.96
$ws.Range("S16").Value = "02/09/2023 21:12"
$ws.Range("T16").Value = 6.56
$ws.Range("U16").Value = "03/09/2023 20:08"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/morocco/botola-pro/berkane-ir-tanger/ILp64n5B/"

# --- Append 6 new match rows (50-55) scraped since the last run ---
# Copy number/style formatting from the last existing data row (49) for the
# "Indice" (A) and "data_partida" (E) columns so the new rows match the
# workbook's existing look (bold/bordered index, datetime-formatted date).
$ws.Range("A49").Copy()
$ws.Range("A50:A55").PasteSpecial(-4122)
$ws.Range("E49").Copy()
$ws.Range("E50:E55").PasteSpecial(-4122)

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "morocco"
$ws.Range("C50").Value = "botola-pro"
$ws.Range("D50").Value = "2023-2024"
$ws.Range("E50").Value = 45227.70833333334
$ws.Range("F50").Value = "Hassania Agadir"
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = "Renaissance Zemamra"
$ws.Range("I50").Value = 2
$ws.Range("J50").Value = 2.32
$ws.Range("K50").Value = "27/10/2023 05:12"
$ws.Range("L50").Value = 3.2
$ws.Range("M50").Value = "28/10/2023 16:58"
$ws.Range("N50").Value = 3.09
$ws.Range("O50").Value = "27/10/2023 05:12"
$ws.Range("P50").Value = 3.03
$ws.Range("Q50").Value = "28/10/2023 16:58"
$ws.Range("R50").Value = 2.86
$ws.Range("S50").Value = "27/10/2023 05:12"
$ws.Range("T50").Value = 2.33
$ws.Range("U50").Value = "28/10/2023 16:58"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/morocco/botola-pro/hassania-agadir-renaissance-zemamra/I103x87I/"

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "morocco"
$ws.Range("C51").Value = "botola-pro"
$ws.Range("D51").Value = "2023-2024"
$ws.Range("E51").Value = 45227.80208333334
$ws.Range("F51").Value = "Olympique de Safi"
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = "Raja Casablanca"
$ws.Range("I51").Value = 2
$ws.Range("J51").Value = 2.61
$ws.Range("K51").Value = "27/10/2023 07:42"
$ws.Range("L51").Value = 5.53
$ws.Range("M51").Value = "28/10/2023 19:10"
$ws.Range("N51").Value = 2.68
$ws.Range("O51").Value = "27/10/2023 07:42"
$ws.Range("P51").Value = 2.87
$ws.Range("Q51").Value = "28/10/2023 19:10"
$ws.Range("R51").Value = 2.8
$ws.Range("S51").Value = "27/10/2023 07:42"
$ws.Range("T51").Value = 1.85
$ws.Range("U51").Value = "28/10/2023 19:10"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/morocco/botola-pro/olympique-de-safi-raja-casablanca/dha7ySMO/"

$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "morocco"
$ws.Range("C52").Value = "botola-pro"
$ws.Range("D52").Value = "2023-2024"
$ws.Range("E52").Value = 45227.89583333334
$ws.Range("F52").Value = "Youssoufia Berrechid"
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = "Union Touarga"
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 2.88
$ws.Range("K52").Value = "28/10/2023 13:43"
$ws.Range("L52").Value = 2.47
$ws.Range("M52").Value = "28/10/2023 21:28"
$ws.Range("N52").Value = 2.96
$ws.Range("O52").Value = "28/10/2023 13:43"
$ws.Range("P52").Value = 2.97
$ws.Range("Q52").Value = "28/10/2023 21:28"
$ws.Range("R52").Value = 2.49
$ws.Range("S52").Value = "28/10/2023 13:43"
$ws.Range("T52").Value = 3.04
$ws.Range("U52").Value = "28/10/2023 21:28"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/morocco/botola-pro/youssoufia-berrechid-union-touarga/Aemk0BaP/"

$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "morocco"
$ws.Range("C53").Value = "botola-pro"
$ws.Range("D53").Value = "2023-2024"
$ws.Range("E53").Value = 45228.66666666666
$ws.Range("F53").Value = "FAR Rabat"
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = "IR Tanger"
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 1.53
$ws.Range("K53").Value = "28/10/2023 21:42"
$ws.Range("L53").Value = 1.43
$ws.Range("M53").Value = "29/10/2023 15:16"
$ws.Range("N53").Value = 4
$ws.Range("O53").Value = "28/10/2023 21:42"
$ws.Range("P53").Value = 4.36
$ws.Range("Q53").Value = "29/10/2023 15:16"
$ws.Range("R53").Value = 5.71
$ws.Range("S53").Value = "28/10/2023 21:42"
$ws.Range("T53").Value = 6.91
$ws.Range("U53").Value = "29/10/2023 15:16"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/morocco/botola-pro/far-rabat-ir-tanger/plvHeV7t/"

$ws.Range("A54").Value = 53
$ws.Range("B54").Value = "morocco"
$ws.Range("C54").Value = "botola-pro"
$ws.Range("D54").Value = "2023-2024"
$ws.Range("E54").Value = 45228.76041666666
$ws.Range("F54").Value = "Moghreb Tetouan"
$ws.Range("G54").Value = 4
$ws.Range("H54").Value = "Jeunesse Sportive Soualem"
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 2.04
$ws.Range("K54").Value = "28/10/2023 18:13"
$ws.Range("L54").Value = 2.02
$ws.Range("M54").Value = "29/10/2023 17:40"
$ws.Range("N54").Value = 2.9
$ws.Range("O54").Value = "28/10/2023 18:13"
$ws.Range("P54").Value = 2.98
$ws.Range("Q54").Value = "29/10/2023 17:46"
$ws.Range("R54").Value = 3.57
$ws.Range("S54").Value = "28/10/2023 18:13"
$ws.Range("T54").Value = 4.09
$ws.Range("U54").Value = "29/10/2023 17:40"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/morocco/botola-pro/moghreb-tetouan-jeunesse-sportive-soualem/SCTLfkNn/"

$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "morocco"
$ws.Range("C55").Value = "botola-pro"
$ws.Range("D55").Value = "2023-2024"
$ws.Range("E55").Value = 45228.85416666666
$ws.Range("F55").Value = "Maghreb Fez"
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = "Berkane"
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3.5
$ws.Range("K55").Value = "28/10/2023 18:13"
$ws.Range("L55").Value = 3.11
$ws.Range("M55").Value = "29/10/2023 20:22"
$ws.Range("N55").Value = 2.89
$ws.Range("O55").Value = "28/10/2023 18:13"
$ws.Range("P55").Value = 2.66
$ws.Range("Q55").Value = "29/10/2023 20:22"
$ws.Range("R55").Value = 2.07
$ws.Range("S55").Value = "28/10/2023 18:13"
$ws.Range("T55").Value = 2.67
$ws.Range("U55").Value = "29/10/2023 20:22"
$ws.Range("V55").Value = "https://www.betexplorer.com/football/morocco/botola-pro/maghreb-fez-berkane/63SPg9xg/"
